$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2018-12-31 00:00:00"
$ws.Range("O2").Value = 50196624.5
$ws.Range("P2").Value = 145.8348089206
$ws.Range("Q2").Value = 718305576.5599999
$ws.Range("R2").Value = 2086.8725247496
$ws.Range("S2").Value = 382656613.39
$ws.Range("T2").Value = 1111.7212492233
$ws.Range("U2").Value = -10602293.6
$ws.Range("V2").Value = -30.8025385507
$ws.Range("W2").Value = 188571.22
$ws.Range("X2").Value = 0.54785054
$ws.Range("Y2").Value = 3825671.57
$ws.Range("Z2").Value = 11.1146135415
$ws.Range("AA2").Value = -5176865.01
$ws.Range("AB2").Value = -15.0401969667
$ws.Range("AC2").Value = 34420194.24
$ws.Range("AD2").Value = 65.5882569708
